$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value updates (Price / Volume(1h) columns) ---
# Numeric-looking Price strings are apostrophe-prefixed so Excel
# stores them as text (matching the original inline-string cells)
# instead of auto-converting them to numbers.
$ws.Range("D2").Value = '43.180.46'
$ws.Range("E2").Value = '  +0.54%  '
$ws.Range("D3").Value = '2.306.54'
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''301.04'
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("D6").Value = '''98.09'
$ws.Range("E6").Value = '  -1.07%  '
$ws.Range("E7").Value = '  +2.96%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '''0.516'
$ws.Range("E9").Value = '  +0.93%  '
$ws.Range("D10").Value = '''36.25'
$ws.Range("E10").Value = '  +0.61%  '
$ws.Range("E11").Value = '  +0.29%  '
$ws.Range("E12").Value = '  +0.71%  '
$ws.Range("D13").Value = '''17.72'
$ws.Range("E13").Value = '  -3.94%  '
$ws.Range("D14").Value = '''6.89'
$ws.Range("E14").Value = '  -0.86%  '
$ws.Range("D15").Value = '2.664.67'
$ws.Range("E15").Value = '  +0.62%  '
$ws.Range("D16").Value = '2.340.22'
$ws.Range("E16").Value = '  +2.16%  '
$ws.Range("D17").Value = '''0.790'
$ws.Range("E17").Value = '  -1.06%  '
$ws.Range("D18").Value = '43.064.20'
$ws.Range("E18").Value = '  +0.54%  '
$ws.Range("D19").Value = '''13.06'
$ws.Range("E19").Value = '  +4.19%  '
$ws.Range("D20").Value = '0.0₃0911'
$ws.Range("E20").Value = '  +0.95%  '
$ws.Range("E21").Value = '  +0.22%  '
$ws.Range("D22").Value = '''68.33'
$ws.Range("E22").Value = '  +0.91%  '
$ws.Range("D23").Value = '''238.17'
$ws.Range("E23").Value = '  +1.14%  '
$ws.Range("D24").Value = '''2.21'
$ws.Range("E24").Value = '  -0.62%  '
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("D26").Value = '''2.43'
$ws.Range("E26").Value = '  -0.72%  '
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("E28").Value = '  +1.32%  '
$ws.Range("D32").Value = '''33.25'
$ws.Range("E32").Value = '  -3.48%  '
$ws.Range("D35").Value = '''18.22'
$ws.Range("E35").Value = '  +3.17%  '
$ws.Range("D36").Value = '''4.73'
$ws.Range("E36").Value = '  +1.53%  '
$ws.Range("D37").Value = '''2.42'
$ws.Range("E37").Value = '  +0.44%  '
$ws.Range("D38").Value = '''0.0696'
$ws.Range("E38").Value = '  +0.83%  '
$ws.Range("D39").Value = '''0.102'
$ws.Range("E39").Value = '  +1.51%  '
$ws.Range("E40").Value = '  +0.10%  '
$ws.Range("E41").Value = '  +1.20%  '
$ws.Range("D42").Value = '''2.76'
$ws.Range("E42").Value = '  -1.80%  '
$ws.Range("D43").Value = '2.016.38'
$ws.Range("E43").Value = '  +1.82%  '
$ws.Range("E44").Value = '  -1.05%  '
$ws.Range("E45").Value = '  -6.99%  '
$ws.Range("D46").Value = '''10.26'
$ws.Range("E46").Value = '  +1.75%  '
$ws.Range("D47").Value = '''17.63'
$ws.Range("E47").Value = '  +0.43%  '
$ws.Range("D48").Value = '''2.85'
$ws.Range("E48").Value = '  -0.89%  '
$ws.Range("D49").Value = '''54.41'
$ws.Range("E49").Value = '  -1.40%  '
$ws.Range("D50").Value = '2.538.00'
$ws.Range("E50").Value = '  +0.91%  '
$ws.Range("E51").Value = '  +0.33%  '

# --- Rows with reordered coins (29-31: Monero/Cosmos/Toncoin -> Cosmos/Toncoin/Monero) ---
# --- and (33-34: Filecoin/FirstDigitalUSD -> FirstDigitalUSD/Filecoin) ---
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").Value = '''9.17'
$ws.Range("E29").Value = '  +0.60%  '

$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '''2.04'
$ws.Range("E30").Value = '  -8.09%  '

$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").Value = '''163.16'
$ws.Range("E31").Value = '  -2.19%  '

$ws.Range("B33").Value = 'FirstDigitalUSD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D33").Value = '''1.00'
$ws.Range("E33").Value = '  +0.00%  '

$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '''5.15'
$ws.Range("E34").Value = '  +2.84%  '

